# The workbook's 16 sheets each hold a statsmodels OLS summary (plain text)
# in cell B2. The summary was regenerated, so the "Date:" and "Time:" lines
# embedded in that text need to be refreshed (all other figures - R-squared,
# coefficients, Log-Likelihood, etc. - stay the same).
#
# Sheet (tab order) -> old Time -> new Time, all sheets move the date from
# "Thu, 02 Jan 2020" to "Sun, 05 Jan 2020".
$wb = $excel.ActiveWorkbook

$oldDate = "Thu, 02 Jan 2020"
$newDate = "Sun, 05 Jan 2020"
$oldTime = "20:48:49"

$newTimes = @{
    1  = "21:22:27"
    2  = "21:22:27"
    3  = "21:22:27"
    4  = "21:22:27"
    5  = "21:22:27"
    6  = "21:22:27"
    7  = "21:22:27"
    8  = "21:22:28"
    9  = "21:22:28"
    10 = "21:22:28"
    11 = "21:22:28"
    12 = "21:22:28"
    13 = "21:22:28"
    14 = "21:22:28"
    15 = "21:22:28"
    16 = "21:22:28"
}

for ($i = 1; $i -le 16; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()
    $newTime = $newTimes[$i]
    $text = $text -replace [regex]::Escape($oldDate), $newDate
    $text = $text -replace [regex]::Escape($oldTime), $newTime
    $cell.Value = $text
}
